# Updated symbol list on Fri Dec 23 05:31:02 UTC 2022 with GitHub Actions
#
# This updates the "Price" (column D) and, for two rows, the "Volume(1h)"
# (column E) values scraped from coinranking.com. The Price cells are
# stored as text (not numbers) in the workbook, so each one is written
# using a leading apostrophe (forces text entry, matching how the value
# already existed) and then the cell style is reset back to "Normal" so
# no stray number-format/quote-prefix style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextPrice($addr, $newValue) {
    $ws.Range($addr).Formula = "'" + $newValue
    $ws.Range($addr).Style = "Normal"
}

Set-TextPrice "D2"  "246.09"
Set-TextPrice "D3"  "22.05"
Set-TextPrice "D4"  "5.429"
Set-TextPrice "D5"  "0.05783"
Set-TextPrice "D6"  "3.389"
Set-TextPrice "D7"  "6.323"
Set-TextPrice "D8"  "0.8183"
Set-TextPrice "D9"  "0.9526"
$ws.Range("E9").Value = "8FTXTokenFTT"
Set-TextPrice "D10" "0.1430"
Set-TextPrice "D11" "0.07498"
Set-TextPrice "D12" "0.03145"
Set-TextPrice "D13" "0.02996"
Set-TextPrice "D14" "4.152"
Set-TextPrice "D15" "0.09406"
Set-TextPrice "D16" "0.001588"
Set-TextPrice "D17" "0.04815"
Set-TextPrice "D18" "0.0005851"
Set-TextPrice "D19" "0.006182"
Set-TextPrice "D20" "0.004123"
Set-TextPrice "D21" "0.0009963"
Set-TextPrice "D23" "3.771"
Set-TextPrice "D24" "2.230"
Set-TextPrice "D27" "0.0004000"
Set-TextPrice "D40" "0.03896"
Set-TextPrice "D41" "0.006341"
Set-TextPrice "D43" "0.003001"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
Set-TextPrice "D44" "0.006690"
Set-TextPrice "D45" "0.00005591"
Set-TextPrice "D47" "0.3801"
